$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569)
    3  = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569)
    4  = @(0.7287194209349384, 1.65323645889881, 3.082599426703578, 6.48142807727062)
    5  = @(0.1554434735375247, 0.05231270169004087, 0.7127328510149897, 6.48142807727062)
    6  = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569)
    7  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569)
    8  = @(0.7287194209349384, 1766.335244827366, 0.1529057820181812, 6.48142807727062)
    9  = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569)
    10 = @(0.1554434735375247, 0.05231270169004087, 3.082599426703578, 0.4998867070740569)
    11 = @(0.02258322285507441, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569)
    12 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569)
    13 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569)
    14 = @(0.1554434735375247, 0.3375848360084654, 3.082599426703578, 0.4998867070740569)
    15 = @(0.1554434735375247, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569)
    16 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $sum = $vals[0] + $vals[1] + $vals[2] + $vals[3]
    $ws.Cells.Item($row, 7).Value = $sum
}
